# Fruta / hortaliza, semanal
# Adds a new week of price data for "Vega Modelo de Temuco - Pomelo" by
# inserting 3 new records at the top of the data block (rows 239-241),
# which pushes all the existing records down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current first data row of this block.
$ws.Rows("239:241").Insert()

# New weekly records to populate into the freshly inserted rows.
$newRows = @(
  @(10, "Vega Modelo de Temuco", "La Araucanía", 44798, 9, "Fruta", 100102, "Cítricos", 100102006, "Pomelo", "Start Ruby", "Especial", 65,  15000, 15000, 15000, "`$/bandeja 15 kilos granel", "Región de O'Higgins", 1000, 15),
  @(10, "Vega Modelo de Temuco", "La Araucanía", 44798, 9, "Fruta", 100102, "Cítricos", 100102006, "Pomelo", "Start Ruby", "Primera",  125, 12000, 12000, 12000, "`$/bandeja 15 kilos granel", "Región de O'Higgins", 800,  15),
  @(10, "Vega Modelo de Temuco", "La Araucanía", 44798, 9, "Fruta", 100102, "Cítricos", 100102006, "Pomelo", "Start Ruby", "Segunda",  85,  8000,  8000,  8000,  "`$/bandeja 15 kilos granel", "Región de O'Higgins", 533,  15)
)

$startRow = 239
for ($i = 0; $i -lt $newRows.Length; $i++) {
  $rowNum = $startRow + $i
  $rowData = $newRows[$i]
  for ($col = 0; $col -lt $rowData.Length; $col++) {
    $ws.Cells.Item($rowNum, $col + 1).Value = $rowData[$col]
  }
}
